$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.987.75'
$ws.Range("E2").Value = '  +5.36%  '
$ws.Range("D3").Value = '2.255.48'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.532'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.72%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +3.33%  '
$ws.Range("E10").Value = '  +9.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '32.34'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0798'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.38%  '
$ws.Range("E13").Value = '  +2.78%  '
$ws.Range("E14").Value = '  +3.57%  '
$ws.Range("D15").Value = '2.603.52'
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("E16").Value = '  +2.17%  '
$ws.Range("D17").Value = '2.290.14'
$ws.Range("E17").Value = '  +2.30%  '
$ws.Range("E18").Value = '  +3.61%  '
$ws.Range("D19").Value = '41.858.88'
$ws.Range("E19").Value = '  +5.18%  '
$ws.Range("E20").Value = '  +9.87%  '
$ws.Range("D21").Value = '0.0₃0902'
$ws.Range("E21").Value = '  +1.91%  '
$ws.Range("E22").Value = '  +3.03%  '
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.48%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +13.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.64'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.53%  '
$ws.Range("E32").Value = '  +6.35%  '
$ws.Range("E34").Value = '  +3.81%  '
$ws.Range("E35").Value = '  +4.19%  '
$ws.Range("E36").Value = '  +2.38%  '
$ws.Range("E37").Value = '  +2.87%  '
$ws.Range("E38").Value = '  +3.73%  '
$ws.Range("E39").Value = '  +4.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.47'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.63%  '
$ws.Range("E41").Value = '  +2.82%  '
$ws.Range("E42").Value = '  +5.74%  '
$ws.Range("D43").Value = '2.055.84'
$ws.Range("E43").Value = '  -2.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.66'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.60%  '
$ws.Range("E45").Value = '  +2.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.90%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.69%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.86'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.20%  '
$ws.Range("E49").Value = '  +3.83%  '
$ws.Range("E50").Value = '  +3.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.34%  '
